$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    # Writing via Formula + Copy/PasteSpecial(values) keeps the numeric-looking
    # literal stored as TEXT (matches the source file's t="str" cells) instead
    # of Excel's normal "looks like a number -> store as Number" coercion.
    $cell = $ws.Range($rangeAddr)
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# New row 3: 1, 2, 3, 4
Set-TextValue "A3" "1"
Set-TextValue "B3" "2"
Set-TextValue "C3" "3"
Set-TextValue "D3" "4"

# Row 4 gains a D4 value of 4 (previously only had A4:C4)
Set-TextValue "D4" "4"

$excel.CutCopyMode = $false

# Remove the leftover duplicate rows 5-12
$ws.Rows("5:12").Delete()

# Sheet view: make sure right-to-left display is off (matches removal of rightToLeft="0")
$ws.DisplayRightToLeft = $false
